# Commentaires front pour la soumission de fichiers excels
#
# - rename the first sheet "m0" -> "m1"
# - add a new, empty sheet "m3" at the very end of the workbook
# - move the selection on the first sheet from K21 to A9
# - keep the first sheet active/selected

$wb = $excel.ActiveWorkbook

# Rename the first worksheet (m0 -> m1)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "m1"

# Add a brand new blank worksheet after the current last sheet (p2),
# so it becomes the 4th / last tab, then name it "m3"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "m3"

# Restore the first sheet as the active one and move the selection to A9
$ws1.Select()
$ws1.Range("A9").Select()
